# "Drop in RMI script results for 3.0"
# Update the BEPEfCT control-lever value (Boolean Exempt Process Emissions
# From Carbon Tax) on the "BEPEfCT" sheet from 1 to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")
$ws.Range("B2").Value = 0
